$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last header cell (G1) onto the new header
# cell (H1) so it shares the exact same style as the other header cells.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Set the new "Save" column header text
$ws.Range("H1").Value = "Save"

# Add the Save values for each data row
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
